# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that are numeric-looking strings (e.g. "1.001") need a leading
# apostrophe so Excel stores them as text (matching the source data's inlineStr type)
# instead of silently converting them to a floating point number.
$textPrefix = [string][char]39

$ws.Range("D2").Value = '30.340.21'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '2.092.70'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("D4").Value = $textPrefix + '1.001'
$ws.Range("E4").Value = '  -0.85%  '
$ws.Range("D5").Value = $textPrefix + '343.56'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("D7").Value = $textPrefix + '0.5240'
$ws.Range("E7").Value = '  +1.67%  '
$ws.Range("D8").Value = $textPrefix + '0.4428'
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = $textPrefix + '54.64'
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("D10").Value = $textPrefix + '0.09339'
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").Value = $textPrefix + '1.170'
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").Value = $textPrefix + '24.88'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = $textPrefix + '8.600'
$ws.Range("E13").Value = '  +3.73%  '
$ws.Range("D14").Value = $textPrefix + '6.904'
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = '2.063.48'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = $textPrefix + '101.31'
$ws.Range("E16").Value = '  +1.95%  '
$ws.Range("D17").Value = $textPrefix + '0.00001161'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").Value = $textPrefix + '1.004'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = $textPrefix + '21.15'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").Value = $textPrefix + '0.06672'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").Value = $textPrefix + '6.345'
$ws.Range("E21").Value = '  +2.44%  '
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '30.308.28'
$ws.Range("E23").Value = '  +1.87%  '
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").Value = $textPrefix + '2.305'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").Value = $textPrefix + '21.82'
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").Value = $textPrefix + '162.84'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = $textPrefix + '2.520'
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("D29").Value = $textPrefix + '133.14'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("D31").Value = $textPrefix + '1.676'
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("D32").Value = $textPrefix + '0.1046'
$ws.Range("D33").Value = $textPrefix + '6.249'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").Value = $textPrefix + '6.789'
$ws.Range("E34").Value = '  +9.58%  '
$ws.Range("D35").Value = $textPrefix + '3.864'
$ws.Range("E35").Value = '  -1.75%  '
$ws.Range("D36").Value = $textPrefix + '10.24'
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = $textPrefix + '0.02636'
$ws.Range("E37").Value = '  +2.45%  '
$ws.Range("D38").Value = $textPrefix + '0.06852'
$ws.Range("E38").Value = '  +2.26%  '
$ws.Range("D39").Value = $textPrefix + '0.6989'
$ws.Range("E39").Value = '  +1.84%  '
$ws.Range("D40").Value = $textPrefix + '1.346'
$ws.Range("E40").Value = '  +3.76%  '
$ws.Range("D41").Value = $textPrefix + '12.54'
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("D42").Value = $textPrefix + '0.2212'
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").Value = $textPrefix + '0.6836'
$ws.Range("E43").Value = '  +3.04%  '
$ws.Range("D44").Value = $textPrefix + '14.35'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = $textPrefix + '2.350'
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("D47").Value = $textPrefix + '1.376'
$ws.Range("E47").Value = '  +18.57%  '
$ws.Range("D48").Value = $textPrefix + '3.638'
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = $textPrefix + '0.00000000342'
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = $textPrefix + '1.209'
$ws.Range("E50").Value = '  +7.95%  '
$ws.Range("D51").Value = $textPrefix + '1.217'
$ws.Range("E51").Value = '  -0.28%  '
